$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New row 5: Roman, Development, 120 min, new remark
$ws.Range("A5").Value = (Get-Date -Year 2014 -Month 5 -Day 5)
$ws.Range("B5").Value = "Roman"
$ws.Range("C5").Value = "Development"
$ws.Range("D5").Value = 120
$ws.Range("E5").Value = "Selecten funktionert komplett serverseitig, mit unterschiedlichen Eingaben"

# New row 6: Roman, Development, 240 min, new remark
$ws.Range("A6").Value = (Get-Date -Year 2014 -Month 5 -Day 6)
$ws.Range("B6").Value = "Roman"
$ws.Range("C6").Value = "Development"
$ws.Range("D6").Value = 240
$ws.Range("E6").Value = "Selecten funktioniert im Zusammenspiel mit der GUI,  Proxy/Client verbessert"

# Remove "applyFont" from the style used by C5/C6 (style index 26) - unbold
$ws.Range("C5:C6").Font.Bold = $false

# Update selection to E6
$ws.Range("E6").Select()

$wb.Save()
